$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jan-2024")

$ws.Range("E2").Value = "Absent"
$ws.Range("C3").Value = "Absent"
$ws.Range("E3").Value = "Present"
$ws.Range("C4").Value = "Absent"
$ws.Range("D4").Value = "Absent"
$ws.Range("E4").Value = "Present"
$ws.Range("E5").Value = "Present"
$ws.Range("B9").Value = "Time : 5PM To 7 PM"

$ws.Activate()
$ws.Range("E5").Select()
